# Adaptations Data Model for Legal info on Resources
# Adds a new "Authorship Resource" column to Table1 (Sheet1) carrying the
# credit line for every character row, and widens/reflows the sheet to
# accommodate it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# --- 1. Add the new table column ------------------------------------------
$newCol = $lo.ListColumns.Add()
$headerCell = $ws.Cells.Item(1, $newCol.Index)
$headerCell.Value = "Authorship Resource"

# --- 2. Fill every data row with the authorship/credit text ---------------
$dataRange = $newCol.DataBodyRange
$dataRange.Value = "Noémi Villars-Amberg, Daniela Subotic"

# --- 3. Match the look of the neighbouring "Keyword" column (border, wrap,
#        vertical-top alignment, text format) by copying its formatting ---
$templateRange = $ws.Range("N2:N55")
$templateRange.Copy()
$dataRange.PasteSpecial(-4122)   # xlPasteFormats

# --- 4. Column layout: un-merge the G width from F, give it its own width -
$ws.Columns.Item(6).ColumnWidth = 135.17
$ws.Columns.Item(7).ColumnWidth = 18.5

# --- 5. A few short rows need to grow so the new column's text fits -------
$ws.Rows.Item(42).RowHeight = 95
$ws.Rows.Item(43).RowHeight = 95
$ws.Rows.Item(50).RowHeight = 95
$ws.Rows.Item(53).RowHeight = 95

# --- 6. Restore the frozen pane (1 column) then leave the new column's
#        data selected, as it was right after the edit ---------------------
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("O2:O55").Select()
